$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'30.098.28"
$ws.Range("E2").Value = "  -1.50%  "
$ws.Range("D3").Value = "'2.107.17"
$ws.Range("E3").Value = "  -0.09%  "
$ws.Range("E4").Value = "  -0.83%  "
$ws.Range("D5").Value = "'349.95"
$ws.Range("E5").Value = "  +4.14%  "
$ws.Range("D7").Value = "'0.5162"
$ws.Range("E7").Value = "  -1.69%  "
$ws.Range("E8").Value = "  -1.38%  "
$ws.Range("D9").Value = "'52.62"
$ws.Range("E9").Value = "  -5.71%  "
$ws.Range("D10").Value = "'0.08973"
$ws.Range("E10").Value = "  -0.51%  "
$ws.Range("D11").Value = "'1.177"
$ws.Range("E11").Value = "  +0.77%  "
$ws.Range("D12").Value = "'25.86"
$ws.Range("E12").Value = "  +5.27%  "
$ws.Range("D13").Value = "'2.113.70"
$ws.Range("E13").Value = "  -0.19%  "
$ws.Range("D14").Value = "'8.265"
$ws.Range("E14").Value = "  +1.87%  "
$ws.Range("D15").Value = "'6.763"
$ws.Range("E15").Value = "  -1.30%  "
$ws.Range("D16").Value = "'99.34"
$ws.Range("E16").Value = "  +2.12%  "
$ws.Range("D17").Value = "'0.00001149"
$ws.Range("E17").Value = "  -2.57%  "
$ws.Range("D18").Value = "'1.004"
$ws.Range("E18").Value = "  -0.75%  "
$ws.Range("D19").Value = "'20.93"
$ws.Range("E19").Value = "  +8.35%  "
$ws.Range("D20").Value = "'0.06678"
$ws.Range("E20").Value = "  -0.15%  "
$ws.Range("D21").Value = "'1.003"
$ws.Range("E21").Value = "  -0.71%  "
$ws.Range("D22").Value = "'6.306"
$ws.Range("E22").Value = "  +0.74%  "
$ws.Range("D23").Value = "'30.185.96"
$ws.Range("E23").Value = "  -1.42%  "
$ws.Range("D25").Value = "'2.356"
$ws.Range("E25").Value = "  -0.22%  "
$ws.Range("D26").Value = "'2.361.41"
$ws.Range("E26").Value = "  -0.14%  "
$ws.Range("D27").Value = "'22.09"
$ws.Range("E27").Value = "  -0.91%  "
$ws.Range("D28").Value = "'2.565"
$ws.Range("E28").Value = "  +2.05%  "
$ws.Range("D29").Value = "'163.19"
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("D30").Value = "'133.78"
$ws.Range("E30").Value = "  +0.19%  "
$ws.Range("D31").Value = "'1.182"
$ws.Range("E31").Value = "  -2.94%  "
$ws.Range("D32").Value = "'0.1068"
$ws.Range("E32").Value = "  -0.03%  "
$ws.Range("D33").Value = "'1.649"
$ws.Range("E33").Value = "  +1.77%  "
$ws.Range("D34").Value = "'6.279"
$ws.Range("E34").Value = "  -1.05%  "
$ws.Range("D35").Value = "'3.979"
$ws.Range("E35").Value = "  +0.31%  "
$ws.Range("D36").Value = "'5.908"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("D37").Value = "'10.23"
$ws.Range("E37").Value = "  -2.41%  "
$ws.Range("D38").Value = "'0.02590"
$ws.Range("E38").Value = "  -0.73%  "
$ws.Range("D39").Value = "'0.06854"
$ws.Range("E39").Value = "  +0.58%  "
$ws.Range("D40").Value = "'0.2327"
$ws.Range("E40").Value = "  +0.68%  "
$ws.Range("D41").Value = "'12.58"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("D42").Value = "'0.6859"
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'1.252"
$ws.Range("E43").Value = "  -0.24%  "
$ws.Range("D44").Value = "'14.25"
$ws.Range("E44").Value = "  +1.67%  "
$ws.Range("D45").Value = "'0.6437"
$ws.Range("E45").Value = "  -0.04%  "
$ws.Range("D46").Value = "'2.306"
$ws.Range("E46").Value = "  +0.06%  "
$ws.Range("D47").Value = "'0.00000000366"
$ws.Range("E47").Value = "  +4.60%  "
$ws.Range("D48").Value = "'3.666"
$ws.Range("E48").Value = "  -0.24%  "
$ws.Range("D49").Value = "'84.20"
$ws.Range("E49").Value = "  +1.30%  "
$ws.Range("E50").Value = "  -2.13%  "
$ws.Range("D51").Value = "'0.07242"
$ws.Range("E51").Value = "  +0.72%  "

# Reset style on D-column cells we touched so the quote-prefix formatting
# introduced by the apostrophe (forcing text entry) does not leave a stray style index
$ws.Range("D2:D3").Style = "Normal"
$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D9:D23").Style = "Normal"
$ws.Range("D25:D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
